# Generate Report for Handback
# Updates the handback status of 9f945c4e-42af-42cd-90b2-35194e932408.md
# from "Ready for handoff" to "Handed back: in sync with en-US" across
# the Overview, zh-cn and de-de sheets, refreshes the handback datetimes,
# and clears the stale error detail message.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-16 00:44:27"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-16 00:44:34"
$wsDeDe.Range("P3").Value = ""
